$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 58 ---
$ws.Cells.Item(58, 4).Value = 44568
$ws.Cells.Item(58, 11).Value = "Flame Seedless"
$ws.Cells.Item(58, 14).Value = 9500
$ws.Cells.Item(58, 15).Value = 10000
$ws.Cells.Item(58, 16).Value = 9750
$ws.Cells.Item(58, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(58, 19).Value = 650
$ws.Cells.Item(58, 20).Value = 15

# --- New row 59: Superior Seedless ---
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44568
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100109
$ws.Cells.Item(59, 8).Value = "Uva"
$ws.Cells.Item(59, 9).Value = 100109001
$ws.Cells.Item(59, 10).Value = "Uva"
$ws.Cells.Item(59, 11).Value = "Superior Seedless"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 360
$ws.Cells.Item(59, 14).Value = 11500
$ws.Cells.Item(59, 15).Value = 12000
$ws.Cells.Item(59, 16).Value = 11750
$ws.Cells.Item(59, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(59, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(59, 19).Value = 783
$ws.Cells.Item(59, 20).Value = 15

# --- New row 60: Red Globe (formerly row 58 data) ---
$ws.Cells.Item(60, 1).Value = 8
$ws.Cells.Item(60, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44357
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100109
$ws.Cells.Item(60, 8).Value = "Uva"
$ws.Cells.Item(60, 9).Value = 100109001
$ws.Cells.Item(60, 10).Value = "Uva"
$ws.Cells.Item(60, 11).Value = "Red Globe"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 400
$ws.Cells.Item(60, 14).Value = 8500
$ws.Cells.Item(60, 15).Value = 9000
$ws.Cells.Item(60, 16).Value = 8750
$ws.Cells.Item(60, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(60, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 19).Value = 486
$ws.Cells.Item(60, 20).Value = 18

# --- New row 61: Red Globe (formerly row 59 data) ---
$ws.Cells.Item(61, 1).Value = 8
$ws.Cells.Item(61, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 44313
$ws.Cells.Item(61, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100109
$ws.Cells.Item(61, 8).Value = "Uva"
$ws.Cells.Item(61, 9).Value = 100109001
$ws.Cells.Item(61, 10).Value = "Uva"
$ws.Cells.Item(61, 11).Value = "Red Globe"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 480
$ws.Cells.Item(61, 14).Value = 6800
$ws.Cells.Item(61, 15).Value = 7000
$ws.Cells.Item(61, 16).Value = 6900
$ws.Cells.Item(61, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(61, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(61, 19).Value = 383
$ws.Cells.Item(61, 20).Value = 18
